$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "崎谷航平"
    3  = "三神佳誠"
    4  = "氏家琉貴"
    5  = "羽賀尚生"
    6  = "島田実"
    7  = "足立耕平"
    8  = "遠藤隼人"
    9  = "富澤天音"
    10 = $null
    11 = "神山修造"
    12 = "志塚惇希"
    13 = "川田涼介"
    14 = "豊島亮"
    15 = "兒島大志郎"
    16 = "山口玲"
    17 = "日高泰聖"
    18 = "白岩詩佑介"
    19 = "Cox Matthew Jonah"
    20 = "Hansen Jakob U"
    21 = "石井海成"
    22 = "Nicholas Tristan Aryasatyo"
    23 = "小溝賢"
    24 = "小野文哉"
    25 = "渡部魁"
    26 = "崎谷航平"
    27 = "三神佳誠"
    28 = "氏家琉貴"
    29 = "羽賀尚生"
    30 = "島田実"
    31 = "足立耕平"
    32 = $null
}

foreach ($row in 2..32) {
    $val = $values[$row]
    if ($null -eq $val) {
        $ws.Cells.Item($row, 2).Value = $null
    } else {
        $ws.Cells.Item($row, 2).Value = $val
    }
}

$ws.Range("B32").Select()
